# Word COM-interop script implementing the LOM3258.docx restructuring:
# paragraph *styles* (Heading2 labels, ListBullet, italics, bold labels, etc.)
# stay exactly where they are; only the underlying text content is shuffled
# between fixed paragraph slots. We therefore rewrite each affected
# paragraph's text in place, scoping every Find/Replace to that paragraph's
# own Range so that identical strings elsewhere in the document can never
# cause a mismatch.

$d = $word.ActiveDocument
$nl = [char]11   # manual line break (<w:br/>) inside a single run

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    $p = $d.Paragraphs.Item($Index)
    $rng = $p.Range
    $ok = $rng.Find.Execute($OldText, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        throw ("Replace failed in paragraph " + $Index + ": '" + $OldText + "'")
    }
}

# --- Objetivos (para 6 / 7): PT+EN objective text moves down to the new
#     "Docente(s) Responsável(eis)" bullet (para 9); the PT+EN "Programa
#     resumido" summary sentences move up into paras 6/7 instead. --------
Replace-InParagraph 6 `
    "Proporcionar ao aluno ingressante de Engenharia Física os conhecimentos práticos de eletrônica e computação física com microcontrolador Arduino visando sua aplicação em projetos científicos e tecnológicos." `
    "Introdução ao Arduino. Conceitos de eletrônica analógica e digital. Montagem de circuitos eletrônicos básicos. Programação e controle de circuitos eletrônicos em linguagem C. Aplicação e desenvolvimento de projetos baseados em Arduino."

Replace-InParagraph 7 `
    "To provide the incoming student of Physical Engineering with practical knowledge of electronics and physical computing with Arduino microcontroller aiming its application in scientific and technological projects" `
    "Introduction to Arduino. Analog and digital electronics concepts. Assembly of basic electronic circuits. Programming and control of electronic circuits in C language. Application and development of projects based on Arduino."

# --- Docente(s) Responsável(eis) bullet (para 9): former teacher-name
#     bullet content is replaced by the PT objective sentence that used to
#     live in para 6. -------------------------------------------------
Replace-InParagraph 9 `
    "519033 - Carlos Yujiro Shigue" `
    "Proporcionar ao aluno ingressante de Engenharia Física os conhecimentos práticos de eletrônica e computação física com microcontrolador Arduino visando sua aplicação em projetos científicos e tecnológicos."

# --- Programa resumido (para 11 / 12): the short PT/EN summary sentences
#     are replaced by the full PT/EN detailed program lists that used to
#     sit under the "Programa" heading. ----------------------------------
$ptPrograma = @(
    "Introdução ao microcontrolador Arduino: histórico, tipos e recursos. Oficina prática: instalação e configuração do IDE Arduino.",
    "Conceitos básicos de eletrônica: funcionamento da protoboard, componentes e instrumentos eletrônicos, medições com multímetro e osciloscópio. Grandezas elétricas: resistência, tensão e corrente. Oficina: montagem de circuitos eletrônicos.",
    "Introdução à linguagem de programação Wiring baseada em C/C++. Tipos de dados, sintaxe básica, controle de fluxo, funções da biblioteca padrão. Principais bibliotecas",
    "Entradas e saídas do Arduino. Sinais analógicos e digitais.",
    "Controle de dispositivos utilizando PWM.",
    "Eletrônica analógica. Conversores analógico-digitais do Arduino. ",
    "Oficina: leitura de dados de sensores. Comunicação serial/USB com o PC. Utilização do Monitor Serial da IDE.",
    "Controle de motor cc e servomotor com PWM. Controle de potência com relé e SSR.",
    "Tópicos avançados: comunicação Ethernet com Arduino. Comunicação sem fio via Bluetooth.",
    "Armazenamento de dados utilizando a EEPROM do ATMega328 e cartão de memória SD.",
    "Desenvolvimento de software de qualidade.",
    "Desenvolvimento de projetos utilizando microcontrolador Arduino."
) -join $nl

Replace-InParagraph 11 `
    "Introdução ao Arduino. Conceitos de eletrônica analógica e digital. Montagem de circuitos eletrônicos básicos. Programação e controle de circuitos eletrônicos em linguagem C. Aplicação e desenvolvimento de projetos baseados em Arduino." `
    $ptPrograma

Replace-InParagraph 12 `
    "Introduction to Arduino. Analog and digital electronics concepts. Assembly of basic electronic circuits. Programming and control of electronic circuits in C language. Application and development of projects based on Arduino." `
    "To provide the incoming student of Physical Engineering with practical knowledge of electronics and physical computing with Arduino microcontroller aiming its application in scientific and technological projects"

# --- Programa (para 14): the big PT detailed-program list is replaced by
#     the short "Aulas expositivas..." sentence that used to be the
#     Avaliação "Método" value. Para 15 (the EN detailed list) stays put. -
$ptProgramaOld = $ptPrograma
Replace-InParagraph 14 `
    $ptProgramaOld `
    "Aulas expositivas, práticas e de realização de projetos."

# --- Avaliação bullet (para 17): three labelled values rotate one slot.
#     Processed from the last value back to the first so that no
#     intermediate text can collide with a not-yet-processed search. ----
$bibliografia = @(
    "BANZI, M. Primeiros passos com o Arduino, São Paulo: O´Reilly Novatec, 2010.",
    "McROBERTS, M. Arduino Básico, São Paulo: Novatec, 2011.",
    "MONK, S. Programação com Arduino, Porto Alegre: Bookman Editora, 2013.",
    "MONK, S. Programação com Arduino II, Porto Alegre: Bookman Editora, 2015.",
    "BLUM, J. Exploring Arduino, New York: John Wiley, 2013."
) -join $nl

Replace-InParagraph 17 `
    "Devido às características da disciplina não será oferecida recuperação." `
    $bibliografia

Replace-InParagraph 17 `
    "Média das notas de trabalhos, atividades e relatório de projeto." `
    "Devido às características da disciplina não será oferecida recuperação."

Replace-InParagraph 17 `
    "Aulas expositivas, práticas e de realização de projetos." `
    "Média das notas de trabalhos, atividades e relatório de projeto."

# --- Bibliografia (para 19): the reference list is replaced by the
#     teacher-name bullet text that used to sit under "Docente(s)
#     Responsável(eis)". --------------------------------------------------
Replace-InParagraph 19 `
    $bibliografia `
    "519033 - Carlos Yujiro Shigue"

Write-Output "Done."
